# Adds UX 1 and UX 2 results
# - Adds a per-row "Average" column (column I) with AVERAGE(B:H) formulas to the
#   "RQ 1", "RQ 2" and "RQ 3" sheets.
# - Relabels a couple of prototype names to disambiguate duplicate "Prototype 4"/
#   "Prototype 5" rows ( [ M ] / [ S ] / [ A ] suffixes ).
# - Cleans up a leftover "NA ( protocol modified )" note, replacing it with "NA ".
# - Updates sheet/window selections so the "RQ 3" tab ends up active.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General")
$wsRQ1 = $wb.Worksheets.Item("RQ 1")
$wsRQ2 = $wb.Worksheets.Item("RQ 2")
$wsRQ3 = $wb.Worksheets.Item("RQ 3")

# ---------------------------------------------------------------------------
# RQ 1 sheet: new "Average" column (I) for the two rating tables, plus
# renamed prototype rows.
# ---------------------------------------------------------------------------
$wsRQ1.Range("I20").Value = "Average"
$wsRQ1.Range("I20").Font.Bold = $true

$wsRQ1.Range("I21").Formula = "=AVERAGE(B21:H21)"
$wsRQ1.Range("I22").Formula = "=AVERAGE(B22:H22)"
$wsRQ1.Range("I23").Formula = "=AVERAGE(B23:H23)"

$wsRQ1.Range("A40").Value = "Prototype 4 [ M ]"
$wsRQ1.Range("A41").Value = "Prototype 5 [ S ]"
$wsRQ1.Range("I40").Formula = "=AVERAGE(B40:H40)"
$wsRQ1.Range("I41").Formula = "=AVERAGE(B41:H41)"

$wsRQ1.Range("A44").Select()

# ---------------------------------------------------------------------------
# RQ 2 sheet: new "Average:" column (I) for the two rating tables, plus
# renamed prototype rows.
# ---------------------------------------------------------------------------
$wsRQ2.Range("I6").Value = "Average:"

$wsRQ2.Range("I7").Formula = "=AVERAGE(B7:H7)"
$wsRQ2.Range("I8").Formula = "=AVERAGE(B8:H8)"
$wsRQ2.Range("I9").Formula = "=AVERAGE(B9:H9)"

$wsRQ2.Range("A17").Value = "Prototype 4 [ A ]"
$wsRQ2.Range("A18").Value = "Prototype 5 [ S ]"
$wsRQ2.Range("I17").Formula = "=AVERAGE(B17:H17)"
$wsRQ2.Range("I18").Formula = "=AVERAGE(B18:H18)"

$wsRQ2.Range("A20").Select()

# ---------------------------------------------------------------------------
# RQ 3 sheet: new "Average:" column (I) for the two rating tables, plus
# cleanup of the "NA ( protocol modified )" note.
# ---------------------------------------------------------------------------
$wsRQ3.Range("I16").Value = "Average:"

$wsRQ3.Range("I17").Formula = "=AVERAGE(B17:H17)"
$wsRQ3.Range("I18").Formula = "=AVERAGE(B18:H18)"

$wsRQ3.Range("E29").Value = "NA "

# ---------------------------------------------------------------------------
# View / selection bookkeeping.
# ---------------------------------------------------------------------------
$wsGeneral.Range("E1:E1048576").Select()

$wsRQ3.Activate()
$wsRQ3.Range("L28").Select()
